# Applies the corrected IFRS figures for 애경유화 (company_list sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 11611
$ws.Range("E2").Value = 412
$ws.Range("F2").Value = 412
$ws.Range("G2").Value = 363
$ws.Range("H2").Value = 262
$ws.Range("I2").Value = 266
$ws.Range("J2").Value = -4
$ws.Range("K2").Value = 3976
$ws.Range("L2").Value = 2162
$ws.Range("M2").Value = 1814
$ws.Range("N2").Value = 1841
$ws.Range("O2").Value = -27
$ws.Range("P2").Value = 160
$ws.Range("Q2").Value = 259
$ws.Range("R2").Value = -52
$ws.Range("S2").Value = -249
$ws.Range("T2").Value = 68
$ws.Range("U2").Value = 191
$ws.Range("V2").Value = 1103
$ws.Range("W2").Value = 3.55
$ws.Range("X2").Value = 2.26
$ws.Range("Y2").Value = 15.3
$ws.Range("Z2").Value = 6.48
$ws.Range("AA2").Value = 119.21
$ws.Range("AB2").Value = 1056.39
$ws.Range("AC2").Value = 830
$ws.Range("AD2").Value = 7.53
$ws.Range("AE2").Value = 5748
$ws.Range("AF2").Value = 1.09
$ws.Range("AG2").Value = 190
$ws.Range("AH2").Value = 3.04
$ws.Range("AI2").Value = 22.89
$ws.Range("AJ2").Value = 32040000

# Row 3
$ws.Range("D3").Value = 9121
$ws.Range("E3").Value = 323
$ws.Range("F3").Value = 323
$ws.Range("G3").Value = 345
$ws.Range("H3").Value = 280
$ws.Range("I3").Value = 283
$ws.Range("J3").Value = -4
$ws.Range("K3").Value = 3327
$ws.Range("L3").Value = 1271
$ws.Range("M3").Value = 2056
$ws.Range("N3").Value = 2087
$ws.Range("O3").Value = -31
$ws.Range("P3").Value = 160
$ws.Range("Q3").Value = 729
$ws.Range("R3").Value = -45
$ws.Range("S3").Value = -673
$ws.Range("T3").Value = 39
$ws.Range("U3").Value = 689
$ws.Range("V3").Value = 495
$ws.Range("W3").Value = 3.54
$ws.Range("X3").Value = 3.06
$ws.Range("Y3").Value = 14.43
$ws.Range("Z3").Value = 7.66
$ws.Range("AA3").Value = 61.85
$ws.Range("AB3").Value = 1208.2
$ws.Range("AC3").Value = 884
$ws.Range("AD3").Value = 6.17
$ws.Range("AE3").Value = 6516
$ws.Range("AF3").Value = 0.84
$ws.Range("AG3").Value = 190
$ws.Range("AH3").Value = 3.48
$ws.Range("AI3").Value = 21.48
$ws.Range("AJ3").Value = 32040000

# Row 4
$ws.Range("D4").Value = 8677
$ws.Range("E4").Value = 750
$ws.Range("F4").Value = 750
$ws.Range("G4").Value = 703
$ws.Range("H4").Value = 534
$ws.Range("I4").Value = 538
$ws.Range("J4").Value = -4
$ws.Range("K4").Value = 4084
$ws.Range("L4").Value = 1553
$ws.Range("M4").Value = 2531
$ws.Range("N4").Value = 2565
$ws.Range("O4").Value = -34
$ws.Range("P4").Value = 160
$ws.Range("Q4").Value = 673
$ws.Range("R4").Value = -235
$ws.Range("S4").Value = -37
$ws.Range("T4").Value = 51
$ws.Range("U4").Value = 623
$ws.Range("V4").Value = 528
$ws.Range("W4").Value = 8.640000000000001
$ws.Range("X4").Value = 6.15
$ws.Range("Y4").Value = 23.11
$ws.Range("Z4").Value = 14.4
$ws.Range("AA4").Value = 61.37
$ws.Range("AB4").Value = 1499.37
$ws.Range("AC4").Value = 1678
$ws.Range("AD4").Value = 7.51
$ws.Range("AE4").Value = 8009
$ws.Range("AF4").Value = 1.57
$ws.Range("AG4").Value = 300
$ws.Range("AH4").Value = 2.38
$ws.Range("AI4").Value = 17.87
$ws.Range("AJ4").Value = 32040000

# Row 5
$ws.Range("D5").Value = 9598
$ws.Range("E5").Value = 751
$ws.Range("F5").Value = 751
$ws.Range("G5").Value = 874
$ws.Range("H5").Value = 647
$ws.Range("I5").Value = 651
$ws.Range("J5").Value = -4
$ws.Range("K5").Value = 4578
$ws.Range("L5").Value = 1572
$ws.Range("M5").Value = 3007
$ws.Range("N5").Value = 3042
$ws.Range("O5").Value = -36
$ws.Range("P5").Value = 160
$ws.Range("Q5").Value = 312
$ws.Range("R5").Value = 33
$ws.Range("S5").Value = -93
$ws.Range("T5").Value = 141
$ws.Range("U5").Value = 172
$ws.Range("V5").Value = 567
$ws.Range("W5").Value = 7.83
$ws.Range("X5").Value = 6.74
$ws.Range("Y5").Value = 23.22
$ws.Range("Z5").Value = 14.95
$ws.Range("AA5").Value = 52.27
$ws.Range("AB5").Value = 1838.24
$ws.Range("AC5").Value = 2032
$ws.Range("AD5").Value = 8.640000000000001
$ws.Range("AE5").Value = 9638
$ws.Range("AF5").Value = 1.82
$ws.Range("AG5").Value = 350
$ws.Range("AH5").Value = 1.99
$ws.Range("AI5").Value = 16.97
$ws.Range("AJ5").Value = 32040000

# Row 6
$ws.Range("D6").Value = 10314
$ws.Range("E6").Value = 520
$ws.Range("F6").Value = 520
$ws.Range("G6").Value = 648
$ws.Range("H6").Value = 527
$ws.Range("I6").Value = 531
$ws.Range("K6").Value = 4985
$ws.Range("L6").Value = 1565
$ws.Range("M6").Value = 3419
$ws.Range("N6").Value = 3450
$ws.Range("P6").Value = 160
$ws.Range("Q6").Value = 464
$ws.Range("R6").Value = -453
$ws.Range("S6").Value = -79
$ws.Range("T6").Value = 360
$ws.Range("U6").Value = 104
$ws.Range("V6").Value = 597
$ws.Range("W6").Value = 5.04
$ws.Range("X6").Value = 5.11
$ws.Range("Y6").Value = 16.36
$ws.Range("Z6").Value = 11.03
$ws.Range("AA6").Value = 45.77
$ws.Range("AB6").Value = 2090.21
$ws.Range("AC6").Value = 1657
$ws.Range("AD6").Value = 4.77
$ws.Range("AE6").Value = 10929
$ws.Range("AF6").Value = 0.72
$ws.Range("AG6").Value = 350
$ws.Range("AH6").Value = 4.43
$ws.Range("AI6").Value = 20.81
$ws.Range("AJ6").Value = 32040000

# Row 7
$ws.Range("D7").Value = 10745
$ws.Range("E7").Value = 491
$ws.Range("G7").Value = 537
$ws.Range("H7").Value = 404
$ws.Range("I7").Value = 441
$ws.Range("K7").Value = 5093
$ws.Range("L7").Value = 1347
$ws.Range("M7").Value = 3746
$ws.Range("P7").Value = 160
$ws.Range("Q7").Value = 347
$ws.Range("R7").Value = -110
$ws.Range("S7").Value = -346
$ws.Range("T7").Value = 83
$ws.Range("U7").Value = 136
$ws.Range("W7").Value = 4.57
$ws.Range("X7").Value = 3.76
$ws.Range("Z7").Value = 8.02
$ws.Range("AA7").Value = 35.96
$ws.Range("AC7").Value = 1376
$ws.Range("AD7").Value = 6.21
$ws.Range("AG7").Value = 350
$ws.Range("AH7").Value = 4.09
$ws.Range("AI7").Value = 25.43
$ws.Range("N7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()

# Row 8
$ws.Range("D8").Value = 11157
$ws.Range("E8").Value = 556
$ws.Range("G8").Value = 612
$ws.Range("H8").Value = 459
$ws.Range("I8").Value = 461
$ws.Range("K8").Value = 5456
$ws.Range("L8").Value = 1360
$ws.Range("M8").Value = 4095
$ws.Range("P8").Value = 160
$ws.Range("Q8").Value = 549
$ws.Range("R8").Value = -107
$ws.Range("S8").Value = -112
$ws.Range("T8").Value = 120
$ws.Range("U8").Value = 429
$ws.Range("W8").Value = 4.98
$ws.Range("X8").Value = 4.11
$ws.Range("Z8").Value = 8.699999999999999
$ws.Range("AA8").Value = 33.21
$ws.Range("AC8").Value = 1439
$ws.Range("AD8").Value = 5.46
$ws.Range("AG8").Value = 350
$ws.Range("AH8").Value = 4.46
$ws.Range("AI8").Value = 24.32
$ws.Range("N8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()

# Row 9
$ws.Range("D9").Value = 11200
$ws.Range("E9").Value = 557
$ws.Range("G9").Value = 613
$ws.Range("H9").Value = 460
$ws.Range("I9").Value = 462
$ws.Range("K9").Value = 5819
$ws.Range("L9").Value = 1373
$ws.Range("M9").Value = 4446
$ws.Range("P9").Value = 160
$ws.Range("Q9").Value = 549
$ws.Range("R9").Value = -109
$ws.Range("S9").Value = -112
$ws.Range("T9").Value = 120
$ws.Range("U9").Value = 429
$ws.Range("W9").Value = 4.97
$ws.Range("X9").Value = 4.11
$ws.Range("Z9").Value = 8.16
$ws.Range("AA9").Value = 30.88
$ws.Range("AC9").Value = 1442
$ws.Range("AD9").Value = 5.44
$ws.Range("AG9").Value = 350
$ws.Range("AH9").Value = 4.46
$ws.Range("AI9").Value = 24.27
$ws.Range("N9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
